$d = $word.ActiveDocument

# 1) Merge "-" + "whats the move app" into a single run.
$d.Content.Find.Execute("-whats the move app", $true, $false, $false, $false, $false, `
    $true, 1, $false, "-whats the move app", 2) | Out-Null

# 2) Merge "-" + "server cloud torrent" into a single run.
$d.Content.Find.Execute("-server cloud torrent", $true, $false, $false, $false, $false, `
    $true, 1, $false, "-server cloud torrent", 2) | Out-Null

# 3) Merge the tab run and the "-for users to be able to create their own accounts" run
#    into one run (tab then text, same run) - use InsertXML on the whole paragraph so
#    the <w:tab/> element is preserved as an element rather than collapsed into <w:t>.
$targetText = "-for users to be able to create their own accounts"
foreach ($para in $d.Paragraphs) {
    $paraText = $para.Range.Text.TrimEnd([char]13)
    if ($paraText -eq "`t" + $targetText) {
        $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
               '<w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:u w:val="single"/></w:rPr></w:pPr>' + `
               '<w:r><w:rPr><w:u w:val="none"/></w:rPr><w:tab/><w:t>' + $targetText + '</w:t></w:r>' + `
               '</w:p>'
        $para.Range.InsertXML($xml) | Out-Null
        break
    }
}

# 4) Remove the trailing empty "-" bullet paragraph (tab + "-") that follows the
#    "Some type of database system..." bullet.
foreach ($para in $d.Paragraphs) {
    $paraText = $para.Range.Text.TrimEnd([char]13)
    if ($paraText -eq "`t-") {
        $para.Range.Delete() | Out-Null
        break
    }
}
